# ---------------------------------------------------------------------------
# Applies the commit "Add files via upload" to Assignment.docx:
#   1. Moves the "_GoBack" bookmark from the end of paragraph 1 to the very
#      end of the document (the run text itself is unaffected - it stays
#      "...". " either way, only the run split + bookmark position change).
#   2. Rewrites paragraph 2 ("However, when doing this solution, ...") to
#      read "I tried to ssh safely to another computer ... accomplish. ".
#   3. Rewrites paragraph 3 ("I tried to ssh safely ...") to read
#      "The solution obviously work but is quite cumbersome when using root
#      lr ssh safely."
#   4. Fills the previously-empty paragraph 4 with the new "redis server"
#      paragraph.
#   5. Appends a brand-new paragraph 5 ("Split-startup works but ...") and
#      re-adds the "_GoBack" bookmark at its very end.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Word package XML namespace wrapper helper for InsertXML -------------------
function New-WordXmlFragment([string]$innerBodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $innerBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Replace the *content* of a whole paragraph (not the paragraph mark) with
# a fresh run structure supplied as raw OOXML <w:r>/<w:proofErr> markup.
function Set-ParagraphRuns($paragraph, [string]$runsXml) {
    $r = $paragraph.Range
    $r.MoveEnd(1, -1) | Out-Null   # exclude the trailing paragraph mark
    $r.Delete()
    $frag = New-WordXmlFragment("<w:p>" + $runsXml + "</w:p>")
    $r.InsertXML($frag)
}

# ---------------------------------------------------------------------------
# Step 1: drop the _GoBack bookmark from its old location (end of paragraph 1)
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Step 2: paragraph 2 -> "I tried to ssh safely ... accomplish. "
# ---------------------------------------------------------------------------
$p2Runs = @'
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">I tried to </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>ssh</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> safely to another computer transferring my public key to the receiving server. Worked as</w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">well. </w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">This is obviously a little bit more tedious to accomplish. </w:t></w:r>
'@
Set-ParagraphRuns $d.Paragraphs.Item(2) $p2Runs

# ---------------------------------------------------------------------------
# Step 3: paragraph 3 -> "The solution obviously work but is quite
#          cumbersome when using root lr ssh safely."
# ---------------------------------------------------------------------------
$p3Runs = @'
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">The solution obviously </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>work</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> but is quite cumbersome when using root </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>lr</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>ssh</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> safely.</w:t></w:r>
'@
Set-ParagraphRuns $d.Paragraphs.Item(3) $p3Runs

# ---------------------------------------------------------------------------
# Step 4: paragraph 4 (previously empty) -> the "redis server" paragraph
# ---------------------------------------------------------------------------
$p4Runs = @'
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">The use of a </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>redis</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> server was easier. I had all the</w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> information </w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">on one server. </w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>If an additional VM gets added it will be easy for other VMs to get the IP-address. The previous solution is not as flexible.</w:t></w:r>
'@
Set-ParagraphRuns $d.Paragraphs.Item(4) $p4Runs

# ---------------------------------------------------------------------------
# Step 5: append a brand new paragraph 5 with the "Split-startup" text, then
#         re-add the _GoBack bookmark collapsed at its very end.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(4).Range.InsertParagraphAfter()

$p5Runs = @'
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Split-</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>startup</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> works but how are the other VMs supposed to talk to each other</w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> with that solution. They can only talk to the proxy, which in turn can talk to the desired VM I guess.</w:t></w:r>
'@
$p5 = $d.Paragraphs.Item(5)
Set-ParagraphRuns $p5 $p5Runs

$endRange = $d.Paragraphs.Item(5).Range.Duplicate
$endRange.MoveEnd(1, -1) | Out-Null
$endRange.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $endRange)

Write-Host "Done."
